$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.995575904846191
$ws.Range("B1").Value = 2.099128246307373
$ws.Range("C1").Value = 7.739345073699951
$ws.Range("D1").Value = 1.001798987388611
$ws.Range("E1").Value = 0.5038214921951294
